$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14: set date value and activity text
$ws.Range("A14").Value = 41133
$ws.Range("B14").Value = "Implemented AMD BitonicSort"

# Update selection to B19 to match the saved view state
$ws.Range("B19").Select()
